$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 7327.375
$ws.Range("J48").Value = 7327.375
$ws.Range("L48").Value = 21982.125
$ws.Range("N48").Value = -22566.125
$ws.Range("H56").Value = 7327.375
$ws.Range("J56").Value = 7327.375
$ws.Range("L56").Value = 21982.125
$ws.Range("N56").Value = -23050.125
$ws.Range("H132").Value = 3253.1875
$ws.Range("I132").Value = 3403.4
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 10210.2
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -7680.200000000001
$ws.Range("N132").Value = -8060
$ws.Range("H137").Value = 807.45
$ws.Range("I137").Value = 751.5
$ws.Range("K137").Value = 2254.5
$ws.Range("M137").Value = 295.5
$ws.Range("H138").Value = 972.30304
$ws.Range("I138").Value = 569.19446
$ws.Range("J138").Value = 2047.2593
$ws.Range("K138").Value = 1707.58338
$ws.Range("L138").Value = 6141.7779
$ws.Range("M138").Value = 3432.41662
$ws.Range("N138").Value = -16421.7779
$ws.Range("H141").Value = 1953.2449
$ws.Range("I141").Value = 634.2439000000001
$ws.Range("K141").Value = 1902.7317
$ws.Range("M141").Value = 3277.2683

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1285.6765
$ws.Range("I2").Value = 1331.2693
$ws.Range("J2").Value = 1137.5
$ws.Range("K2").Value = 1331.2693
$ws.Range("L2").Value = 1137.5
$ws.Range("M2").Value = -1218.2693
$ws.Range("N2").Value = -1363.5
$ws.Range("H32").Value = 25690.328
$ws.Range("I32").Value = 28829.367
$ws.Range("J32").Value = 19281.459
$ws.Range("K32").Value = 28829.367
$ws.Range("L32").Value = 19281.459
$ws.Range("M32").Value = -28542.367
$ws.Range("N32").Value = -19855.459
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H116").Value = 1285.6765
$ws.Range("I116").Value = 1331.2693
$ws.Range("J116").Value = 1137.5
$ws.Range("K116").Value = 1331.2693
$ws.Range("L116").Value = 1137.5
$ws.Range("M116").Value = 962.7307000000001
$ws.Range("N116").Value = -5725.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1285.6765
$ws.Range("I3").Value = 1331.2693
$ws.Range("J3").Value = 1137.5
$ws.Range("K3").Value = 1331.2693
$ws.Range("L3").Value = 1137.5
$ws.Range("M3").Value = -1217.2693
$ws.Range("N3").Value = -1365.5
$ws.Range("H97").Value = 17607
$ws.Range("I97").Value = 1214
$ws.Range("J97").Value = 34000
$ws.Range("K97").Value = 1214
$ws.Range("L97").Value = 34000
$ws.Range("M97").Value = -223
$ws.Range("N97").Value = -35982
$ws.Range("H134").Value = 19072.293
$ws.Range("I134").Value = 1422.0435
$ws.Range("J134").Value = 86731.586
$ws.Range("K134").Value = 4266.1305
$ws.Range("L134").Value = 260194.758
$ws.Range("M134").Value = -1731.1305
$ws.Range("N134").Value = -265264.758

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2205.869
$ws.Range("I31").Value = 2162.3865
$ws.Range("J31").Value = 2318.4119
$ws.Range("K31").Value = 2162.3865
$ws.Range("L31").Value = 2318.4119
$ws.Range("M31").Value = -1867.3865
$ws.Range("N31").Value = -2908.4119
$ws.Range("H34").Value = 2205.869
$ws.Range("I34").Value = 2162.3865
$ws.Range("J34").Value = 2318.4119
$ws.Range("K34").Value = 2162.3865
$ws.Range("L34").Value = 2318.4119
$ws.Range("M34").Value = -1960.3865
$ws.Range("N34").Value = -2722.4119
$ws.Range("H132").Value = 1414.6034
$ws.Range("I132").Value = 997.8095
$ws.Range("J132").Value = 2508.6875
$ws.Range("K132").Value = 2993.4285
$ws.Range("L132").Value = 7526.0625
$ws.Range("M132").Value = -463.4285
$ws.Range("N132").Value = -12586.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2941286.8
$ws.Range("I2").Value = 6536082
$ws.Range("J2").Value = 90.545456
$ws.Range("K2").Value = 39216492
$ws.Range("L2").Value = 543.272736
$ws.Range("M2").Value = -39216379
$ws.Range("N2").Value = -769.272736
$ws.Range("H5").Value = 765.61536
$ws.Range("I5").Value = 431.75
$ws.Range("J5").Value = 1299.8
$ws.Range("K5").Value = 1295.25
$ws.Range("L5").Value = 3899.4
$ws.Range("M5").Value = -1183.25
$ws.Range("N5").Value = -4123.4
$ws.Range("H14").Value = 126.17647
$ws.Range("I14").Value = 126.17647
$ws.Range("K14").Value = 378.52941
$ws.Range("M14").Value = -205.52941
$ws.Range("H63").Value = 4228.364
$ws.Range("J63").Value = 4955.5557
$ws.Range("L63").Value = 14866.6671
$ws.Range("N63").Value = -16364.6671
$ws.Range("H64").Value = 2500
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 2800
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 8400
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = -8940
$ws.Range("H66").Value = 4228.364
$ws.Range("J66").Value = 4955.5557
$ws.Range("L66").Value = 44600.0013
$ws.Range("N66").Value = -52088.0013
$ws.Range("H67").Value = 2500
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 2800
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 8400
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -10272
$ws.Range("H93").Value = 3000
$ws.Range("J93").Value = 3000
$ws.Range("L93").Value = 9000
$ws.Range("N93").Value = -12744
$ws.Range("H114").Value = 2108.75
$ws.Range("I114").Value = 164
$ws.Range("J114").Value = 2757
$ws.Range("K114").Value = 492
$ws.Range("L114").Value = 8271
$ws.Range("M114").Value = 2762
$ws.Range("N114").Value = -14779
$ws.Range("H129").Value = 50070.24
$ws.Range("I129").Value = 1393.3334
$ws.Range("J129").Value = 69541
$ws.Range("K129").Value = 4180.0002
$ws.Range("L129").Value = 208623
$ws.Range("M129").Value = 819.9997999999996
$ws.Range("N129").Value = -218623
$ws.Range("H131").Value = 7173036.5
$ws.Range("J131").Value = 17189.89
$ws.Range("L131").Value = 51569.67
$ws.Range("N131").Value = -61649.67
$ws.Range("H135").Value = 765.61536
$ws.Range("I135").Value = 431.75
$ws.Range("J135").Value = 1299.8
$ws.Range("K135").Value = 3885.75
$ws.Range("L135").Value = 11698.2
$ws.Range("M135").Value = -1350.75
$ws.Range("N135").Value = -16768.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4459.3125
$ws.Range("I70").Value = 4034.2856
$ws.Range("K70").Value = 4034.2856
$ws.Range("M70").Value = -3764.2856
$ws.Range("H73").Value = 4459.3125
$ws.Range("I73").Value = 4034.2856
$ws.Range("K73").Value = 4034.2856
$ws.Range("M73").Value = -3098.2856
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H132").Value = 1899.6333
$ws.Range("I132").Value = 1771.2084
$ws.Range("K132").Value = 5313.6252
$ws.Range("M132").Value = -2783.6252

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 18564
$ws.Range("J80").Value = 18564
$ws.Range("L80").Value = 18564
$ws.Range("N80").Value = -20810
$ws.Range("H83").Value = 18564
$ws.Range("J83").Value = 18564
$ws.Range("L83").Value = 55692
$ws.Range("N83").Value = -66924
$ws.Range("H92").Value = 24908.908
$ws.Range("J92").Value = 24908.908
$ws.Range("L92").Value = 24908.908
$ws.Range("N92").Value = -29900.908
$ws.Range("H122").Value = 7916.353
$ws.Range("I122").Value = 9582.923000000001
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 28748.769
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -26298.769
$ws.Range("N122").Value = -12400
